$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 16 with the "Merge Two Sorted Lists" entry ---
# Seed row 16 with row 15's formatting (fill/border/alignment) first so the
# new cells pick up the same category styling as the rest of the table.
$ws.Range("A15:G15").Copy()
$ws.Range("A16:G16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write the cell values (order mirrors the authoring order so shared-string
# indices land url, highlight, question - same as the source commit).
$ws.Range("D16").Value = "https://rebrand.ly/4t4p9zb"
$ws.Range("G16").Value = "Use two pointers on each list and keep comparing until one goes invalid. Append the remaining to the result."
$ws.Range("C16").Value = "Merge Two Sorted Lists"
$ws.Range("A16").Value = 21
$ws.Range("B16").Value = "Easy"
$ws.Range("E16").Value = "Pointers"
$ws.Range("F16").Value = "O(n)"

# Hyperlink the url cell, then restore its format (Hyperlinks.Add stomps the
# cell style with a generic hyperlink style) back to match column D's look.
$ws.Hyperlinks.Add($ws.Range("D16"), "https://rebrand.ly/4t4p9zb")
$ws.Range("D15").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Widen column G to fit the new, longer highlight text ---
$ws.Columns("G").ColumnWidth = 101

# --- Match the saved selection state from the diff ---
[void]$ws.Range("D19").Select()
